$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch calculation mode back to automatic (removes calcMode="manual")
$excel.Calculation = -4105  # xlCalculationAutomatic

# Update data values that changed from the refreshed query (simulating a
# "Refresh All" of the external data connection table).
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 4
$ws.Range("D6").Value = 5
$ws.Range("F6").Value = 3
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 5

# Column widths widened slightly (as would happen on a table/data refresh's
# auto-fit of the column headers/content). Values chosen so the engine's
# column-width quantization lands on the closest representable width to the
# target stored width.
$ws.Columns.Item(1).ColumnWidth = 11.833333333333334
$ws.Columns.Item(2).ColumnWidth = 43.666666666666664
$ws.Columns.Item(3).ColumnWidth = 15.5
$ws.Columns.Item(4).ColumnWidth = 11.5
$ws.Columns.Item(5).ColumnWidth = 15.833333333333334
$ws.Columns.Item(6).ColumnWidth = 15.166666666666666

# Move the active selection to D25
$ws.Range("D25").Select()
